$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.576.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.48%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.645.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.40%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '324.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.523'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.94%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.547'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.21'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.08'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0815'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.08%  '
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.052.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.651.72'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.865'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.526.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0952'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '273.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.59'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.33%  '
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.138'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.70'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0798'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.05'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.65%  '
$ws.Range("E41").Value = '  +0.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.30'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.31%  '
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0315'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.101.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.75%  '
$ws.Range("E48").Value = '  +5.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.92'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.56%  '
